$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh6 = $s.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$seg = $tr6.Characters(29, 5)
Write-Host "Seg before insert: [$($seg.Text)]"
$seg.InsertAfter("[]")
Write-Host "Whole after: [$($tr6.Text)]"
